$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2734.1177
$ws.Range("J17").Value = 2734.1177
$ws.Range("L17").Value = 8202.3531
$ws.Range("N17").Value = -8538.3531

$ws.Range("H18").Value = 837
$ws.Range("I18").Value = 837
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 837
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -553
$ws.Range("N18").ClearContents()

$ws.Range("H33").Value = 141.8
$ws.Range("J33").Value = 174.75
$ws.Range("L33").Value = 174.75
$ws.Range("N33").Value = -632.75

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H112").Value = 2825
$ws.Range("I112").Value = 1350
$ws.Range("J112").Value = 3193.75
$ws.Range("K112").Value = 4050
$ws.Range("L112").Value = 9581.25
$ws.Range("M112").Value = -2942
$ws.Range("N112").Value = -11797.25

$ws.Range("H115").Value = 359
$ws.Range("I115").Value = 359
$ws.Range("K115").Value = 1077
$ws.Range("M115").Value = 490

$ws.Range("H129").Value = 1259.5714
$ws.Range("J129").Value = 2675.3333
$ws.Range("L129").Value = 8025.999899999999
$ws.Range("N129").Value = -18025.9999

$ws.Range("H138").Value = 2819.3572
$ws.Range("I138").Value = 2158.6667
$ws.Range("K138").Value = 6476.000100000001
$ws.Range("M138").Value = -1336.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -850
$ws.Range("N30").ClearContents()

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H88").Value = 2730.625
$ws.Range("J88").Value = 3066.1538
$ws.Range("L88").Value = 3066.1538
$ws.Range("N88").Value = -3878.1538

$ws.Range("H91").Value = 2730.625
$ws.Range("J91").Value = 3066.1538
$ws.Range("L91").Value = 3066.1538
$ws.Range("N91").Value = -5874.1538

$ws.Range("H110").Value = 1665.6666
$ws.Range("I110").Value = 1665.6666
$ws.Range("K110").Value = 1665.6666
$ws.Range("M110").Value = 379.3334

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5130.9
$ws.Range("I99").Value = 5716.6665
$ws.Range("J99").Value = 4252.25
$ws.Range("K99").Value = 5716.6665
$ws.Range("L99").Value = 4252.25
$ws.Range("M99").Value = -4218.6665
$ws.Range("N99").Value = -7248.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 28347.166
$ws.Range("J50").Value = 29800
$ws.Range("L50").Value = 29800
$ws.Range("N50").Value = -31050

$ws.Range("H51").Value = 22681.666
$ws.Range("I51").Value = 12090
$ws.Range("J51").Value = 24800
$ws.Range("K51").Value = 12090
$ws.Range("L51").Value = 24800
$ws.Range("M51").Value = -11354
$ws.Range("N51").Value = -26272

$ws.Range("H60").Value = 22182.166
$ws.Range("I60").Value = 10093
$ws.Range("J60").Value = 24600
$ws.Range("K60").Value = 10093
$ws.Range("L60").Value = 24600
$ws.Range("M60").Value = -9582
$ws.Range("N60").Value = -25622

$ws.Range("H61").Value = 22681.666
$ws.Range("I61").Value = 12090
$ws.Range("J61").Value = 24800
$ws.Range("K61").Value = 12090
$ws.Range("L61").Value = 24800
$ws.Range("M61").Value = -11742
$ws.Range("N61").Value = -25496

$ws.Range("H86").Value = 10178.889
$ws.Range("I86").Value = 9370.333000000001
$ws.Range("K86").Value = 9370.333000000001
$ws.Range("M86").Value = -8247.333000000001

$ws.Range("H89").Value = 10178.889
$ws.Range("I89").Value = 9370.333000000001
$ws.Range("K89").Value = 46851.665
$ws.Range("M89").Value = -41235.665

$ws.Range("H105").Value = 1500
$ws.Range("I105").Value = 1250
$ws.Range("K105").Value = 1250
$ws.Range("M105").Value = 497

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2233.3333
$ws.Range("I4").Value = 2161.5386
$ws.Range("J4").Value = 2420
$ws.Range("K4").Value = 6484.6158
$ws.Range("L4").Value = 7260
$ws.Range("M4").Value = -6372.6158
$ws.Range("N4").Value = -7484

$ws.Range("H11").Value = 286573.72
$ws.Range("I11").Value = 500752.5
$ws.Range("J11").Value = 1002
$ws.Range("K11").Value = 1502257.5
$ws.Range("L11").Value = 3006
$ws.Range("M11").Value = -1502117.5
$ws.Range("N11").Value = -3286

$ws.Range("H16").Value = 626.6667
$ws.Range("I16").Value = 470
$ws.Range("J16").Value = 940
$ws.Range("K16").Value = 1410
$ws.Range("L16").Value = 2820
$ws.Range("M16").Value = -1237
$ws.Range("N16").Value = -3166

$ws.Range("H17").Value = 870.5714
$ws.Range("I17").Value = 48
$ws.Range("J17").Value = 1199.6
$ws.Range("K17").Value = 144
$ws.Range("L17").Value = 3598.8
$ws.Range("M17").Value = 25
$ws.Range("N17").Value = -3936.8

$ws.Range("H39").Value = 1275
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H49").Value = 941.2
$ws.Range("I49").Value = 368.66666
$ws.Range("J49").Value = 1800
$ws.Range("K49").Value = 1105.99998
$ws.Range("L49").Value = 5400
$ws.Range("M49").Value = -949.9999800000001
$ws.Range("N49").Value = -5712

$ws.Range("H55").Value = 2631.3333
$ws.Range("I55").Value = 995
$ws.Range("J55").Value = 3449.5
$ws.Range("K55").Value = 2985
$ws.Range("L55").Value = 10348.5
$ws.Range("M55").Value = -2808
$ws.Range("N55").Value = -10702.5

$ws.Range("H115").Value = 212
$ws.Range("I115").Value = 212
$ws.Range("K115").Value = 636
$ws.Range("M115").Value = 539

$ws.Range("H131").Value = 810.8333
$ws.Range("I131").Value = 810.8333
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2432.4999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2607.5001
$ws.Range("N131").ClearContents()

$ws.Range("H140").Value = 1086.6666
$ws.Range("I140").Value = 1086.6666
$ws.Range("K140").Value = 3259.9998
$ws.Range("M140").Value = 1920.0002

$ws.Range("H141").Value = 5000
$ws.Range("I141").Value = 5000
$ws.Range("K141").Value = 15000
$ws.Range("M141").Value = -9820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1425
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 1425
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H40").Value = 1913.5714
$ws.Range("I40").Value = 1913.5714
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1913.5714
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1777.5714
$ws.Range("N40").ClearContents()

$ws.Range("H55").Value = 284.66666
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 308.85715
$ws.Range("K55").Value = 200
$ws.Range("L55").Value = 308.85715
$ws.Range("M55").Value = -27
$ws.Range("N55").Value = -654.85715

$ws.Range("H82").Value = 2196
$ws.Range("I82").Value = 1894
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 1894
$ws.Range("L82").Value = 2800
$ws.Range("M82").Value = -1533
$ws.Range("N82").Value = -3522

$ws.Range("H85").Value = 2196
$ws.Range("I85").Value = 1894
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 1894
$ws.Range("L85").Value = 2800
$ws.Range("M85").Value = -646
$ws.Range("N85").Value = -5296

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H122").Value = 3490
$ws.Range("I122").Value = 3490
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10470
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8020
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 966.6667
$ws.Range("I126").Value = 966.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2900.0001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -430.0001000000002
$ws.Range("N126").ClearContents()
